# Auto-generated edit script: updates market-value columns (H-N) across all 8 sheets
# per the commit diff. Values are set directly; cells that should be removed entirely
# (no corresponding <c> element) are cleared with ClearContents().

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 400
$ws.Range("J12").Value = 0
$ws.Range("M12").Value = -230
$ws.Range("I12").Value = 400
$ws.Range("N12").ClearContents()
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 0
$ws.Range("I21").Value = 9800
$ws.Range("K21").Value = 9800
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -9332
$ws.Range("J21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H21").Value = 9800
$ws.Range("I23").Value = 9800
$ws.Range("N23").ClearContents()
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 9800
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -9566
$ws.Range("H23").Value = 9800
$ws.Range("N40").Value = -7799.5
$ws.Range("L40").Value = 7449.5
$ws.Range("M40").Value = -5824
$ws.Range("I40").Value = 5999
$ws.Range("H40").Value = 6966
$ws.Range("K40").Value = 5999
$ws.Range("J40").Value = 7449.5
$ws.Range("I51").Value = 15333.333
$ws.Range("J51").Value = 10333.333
$ws.Range("M51").Value = -14849.333
$ws.Range("L51").Value = 10333.333
$ws.Range("N51").Value = -11301.333
$ws.Range("K51").Value = 15333.333
$ws.Range("H51").Value = 12833.333
$ws.Range("H58").Value = 1716.75
$ws.Range("I58").Value = 346.8
$ws.Range("M58").Value = -890.4000000000001
$ws.Range("K58").Value = 1040.4
$ws.Range("M70").Value = -6629.25
$ws.Range("K70").Value = 6899.25
$ws.Range("H70").Value = 2156.2173
$ws.Range("J70").Value = 2126
$ws.Range("L70").Value = 6378
$ws.Range("I70").Value = 2299.75
$ws.Range("N70").Value = -6918
$ws.Range("H73").Value = 2156.2173
$ws.Range("M73").Value = -5963.25
$ws.Range("I73").Value = 2299.75
$ws.Range("J73").Value = 2126
$ws.Range("L73").Value = 6378
$ws.Range("K73").Value = 6899.25
$ws.Range("N73").Value = -8250
$ws.Range("M76").Value = -5079.2
$ws.Range("H76").Value = 6786.6665
$ws.Range("K76").Value = 5394.2
$ws.Range("I76").Value = 5394.2
$ws.Range("M79").Value = -4302.2
$ws.Range("K79").Value = 5394.2
$ws.Range("H79").Value = 6786.6665
$ws.Range("I79").Value = 5394.2
$ws.Range("H82").Value = 6672.5713
$ws.Range("I82").Value = 4284.8335
$ws.Range("M82").Value = -12448.5005
$ws.Range("K82").Value = 12854.5005
$ws.Range("I85").Value = 4284.8335
$ws.Range("K85").Value = 12854.5005
$ws.Range("H85").Value = 6672.5713
$ws.Range("M85").Value = -11450.5005
$ws.Range("I92").Value = 100158.3
$ws.Range("M92").Value = -98910.3
$ws.Range("K92").Value = 100158.3
$ws.Range("H92").Value = 77446.62
$ws.Range("I94").Value = 8891.333000000001
$ws.Range("H94").Value = 8891.333000000001
$ws.Range("M94").Value = -8440.333000000001
$ws.Range("K94").Value = 8891.333000000001
$ws.Range("I106").Value = 3917.0833
$ws.Range("K106").Value = 3917.0833
$ws.Range("H106").Value = 4275.385
$ws.Range("M106").Value = -3286.0833
$ws.Range("I107").Value = 926
$ws.Range("M107").Value = 994
$ws.Range("H107").Value = 1540.8
$ws.Range("K107").Value = 926
$ws.Range("N111").Value = -7604
$ws.Range("I111").Value = 500
$ws.Range("J111").Value = 490
$ws.Range("L111").Value = 1470
$ws.Range("K111").Value = 1500
$ws.Range("H111").Value = 495
$ws.Range("M111").Value = 1567
$ws.Range("N116").Value = -12084
$ws.Range("J116").Value = 5200
$ws.Range("H116").Value = 5291.7144
$ws.Range("L116").Value = 5200
$ws.Range("H125").Value = 4199
$ws.Range("M125").Value = -16971
$ws.Range("K125").Value = 19431
$ws.Range("I125").Value = 2159
$ws.Range("M137").Value = -101077.452
$ws.Range("I137").Value = 34542.484
$ws.Range("H137").Value = 28235.77
$ws.Range("K137").Value = 103627.452
$ws.Range("H138").Value = 2487.4
$ws.Range("I138").Value = 2197.6428
$ws.Range("N138").Value = -18643.667
$ws.Range("L138").Value = 8363.667000000001
$ws.Range("J138").Value = 2787.889
$ws.Range("K138").Value = 6592.928400000001
$ws.Range("M138").Value = -1452.928400000001
$ws.Range("H141").Value = 862.1724
$ws.Range("I141").Value = 862.1724
$ws.Range("M141").Value = 2593.4828
$ws.Range("K141").Value = 2586.5172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 30245.568
$ws.Range("K32").Value = 30245.568
$ws.Range("H32").Value = 26660.936
$ws.Range("M32").Value = -29958.568
$ws.Range("H45").Value = 3876.1428
$ws.Range("J45").Value = 4391.8213
$ws.Range("N45").Value = -5145.8213
$ws.Range("L45").Value = 4391.8213
$ws.Range("M61").Value = -8691.700000000001
$ws.Range("K61").Value = 8903.700000000001
$ws.Range("I61").Value = 8903.700000000001
$ws.Range("H61").Value = 8717.809999999999
$ws.Range("H74").Value = 1995.9656
$ws.Range("N74").Value = -3847.8
$ws.Range("K74").Value = 1974.3334
$ws.Range("I74").Value = 1974.3334
$ws.Range("L74").Value = 2099.8
$ws.Range("M74").Value = -1100.3334
$ws.Range("J74").Value = 2099.8
$ws.Range("L77").Value = 10499
$ws.Range("J77").Value = 2099.8
$ws.Range("N77").Value = -19235
$ws.Range("I77").Value = 1974.3334
$ws.Range("K77").Value = 9871.666999999999
$ws.Range("H77").Value = 1995.9656
$ws.Range("M77").Value = -5503.666999999999
$ws.Range("M88").Value = -2293
$ws.Range("I88").Value = 2699
$ws.Range("N88").Value = -4497.5
$ws.Range("K88").Value = 2699
$ws.Range("J88").Value = 3685.5
$ws.Range("H88").Value = 3595.818
$ws.Range("L88").Value = 3685.5
$ws.Range("H91").Value = 3595.818
$ws.Range("K91").Value = 2699
$ws.Range("I91").Value = 2699
$ws.Range("J91").Value = 3685.5
$ws.Range("L91").Value = 3685.5
$ws.Range("M91").Value = -1295
$ws.Range("N91").Value = -6493.5
$ws.Range("K132").Value = 101891.625
$ws.Range("H132").Value = 25455.25
$ws.Range("I132").Value = 33963.875
$ws.Range("M132").Value = -99361.625
$ws.Range("I136").Value = 8903.700000000001
$ws.Range("M136").Value = -24161.1
$ws.Range("H136").Value = 8717.809999999999
$ws.Range("K136").Value = 26711.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M20").Value = -4000717.8
$ws.Range("N20").Value = -3200.8572
$ws.Range("K20").Value = 4000964.8
$ws.Range("J20").Value = 2706.8572
$ws.Range("I20").Value = 4000964.8
$ws.Range("L20").Value = 2706.8572
$ws.Range("H20").Value = 3126346
$ws.Range("I22").Value = 143984.58
$ws.Range("H22").Value = 126174
$ws.Range("K22").Value = 143984.58
$ws.Range("M22").Value = -143811.58
$ws.Range("M86").Value = -1897.8462
$ws.Range("K86").Value = 3020.8462
$ws.Range("I86").Value = 3020.8462
$ws.Range("N86").Value = -5553.4
$ws.Range("J86").Value = 3307.4
$ws.Range("L86").Value = 3307.4
$ws.Range("H86").Value = 3100.4443
$ws.Range("J89").Value = 3307.4
$ws.Range("K89").Value = 15104.231
$ws.Range("N89").Value = -27769
$ws.Range("H89").Value = 3100.4443
$ws.Range("L89").Value = 16537
$ws.Range("I89").Value = 3020.8462
$ws.Range("M89").Value = -9488.231
$ws.Range("I94").Value = 1682
$ws.Range("H94").Value = 1891.7906
$ws.Range("M94").Value = -1231
$ws.Range("K94").Value = 1682
$ws.Range("I107").Value = 2892.8462
$ws.Range("M107").Value = -972.8462
$ws.Range("H107").Value = 3252.3125
$ws.Range("K107").Value = 2892.8462
$ws.Range("M134").Value = -6572.7999
$ws.Range("K134").Value = 9107.7999
$ws.Range("H134").Value = 3146.2122
$ws.Range("I134").Value = 3035.9333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 6127.25
$ws.Range("J12").Value = 3750
$ws.Range("M12").Value = -8334.5
$ws.Range("I12").Value = 8504.5
$ws.Range("N12").Value = -4090
$ws.Range("K12").Value = 8504.5
$ws.Range("L12").Value = 3750
$ws.Range("J31").Value = 4127.857
$ws.Range("K31").Value = 1126.2106
$ws.Range("H31").Value = 1934.3462
$ws.Range("N31").Value = -4717.857
$ws.Range("I31").Value = 1126.2106
$ws.Range("L31").Value = 4127.857
$ws.Range("M31").Value = -831.2106000000001
$ws.Range("I34").Value = 1126.2106
$ws.Range("H34").Value = 1934.3462
$ws.Range("M34").Value = -924.2106000000001
$ws.Range("L34").Value = 4127.857
$ws.Range("J34").Value = 4127.857
$ws.Range("N34").Value = -4531.857
$ws.Range("K34").Value = 1126.2106
$ws.Range("H58").Value = 52511.9
$ws.Range("I58").Value = 73589
$ws.Range("M58").Value = -73386
$ws.Range("K58").Value = 73589
$ws.Range("M70").ClearContents()
$ws.Range("K70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("I107").Value = 691.8461
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1228.1539
$ws.Range("H107").Value = 691.8461
$ws.Range("N107").ClearContents()
$ws.Range("K107").Value = 691.8461
$ws.Range("I136").Value = 73589
$ws.Range("M136").Value = -218217
$ws.Range("H136").Value = 52511.9
$ws.Range("K136").Value = 220767

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L2").Value = 158.000004
$ws.Range("J2").Value = 26.333334
$ws.Range("K2").Value = 402
$ws.Range("N2").Value = -384.000004
$ws.Range("M2").Value = -289
$ws.Range("H2").Value = 51.75
$ws.Range("I2").Value = 67
$ws.Range("H12").Value = 147.61905
$ws.Range("J12").Value = 148.35715
$ws.Range("M12").Value = -265.42855
$ws.Range("I12").Value = 146.14285
$ws.Range("N12").Value = -791.0714499999999
$ws.Range("K12").Value = 438.42855
$ws.Range("L12").Value = 445.07145
$ws.Range("I23").Value = 300
$ws.Range("N23").ClearContents()
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 900
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -665
$ws.Range("H23").Value = 300
$ws.Range("I34").Value = 995
$ws.Range("H34").Value = 965
$ws.Range("M34").Value = -2901
$ws.Range("K34").Value = 2985
$ws.Range("H39").Value = 6746.625
$ws.Range("N39").Value = -20827.875
$ws.Range("L39").Value = 20239.875
$ws.Range("J39").Value = 6746.625
$ws.Range("N55").Value = -6215.571599999999
$ws.Range("L55").Value = 5861.571599999999
$ws.Range("J55").Value = 1953.8572
$ws.Range("H55").Value = 2069.9473
$ws.Range("L114").Value = 3273
$ws.Range("H114").Value = 861.2308
$ws.Range("J114").Value = 1091
$ws.Range("N114").Value = -9781
$ws.Range("K119").Value = 5992.9998
$ws.Range("H119").Value = 10998.167
$ws.Range("I119").Value = 1997.6666
$ws.Range("M119").Value = -1154.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L63").Value = 58332.332
$ws.Range("H63").Value = 58332.332
$ws.Range("J63").Value = 58332.332
$ws.Range("N63").Value = -59704.332
$ws.Range("L66").Value = 174996.996
$ws.Range("H66").Value = 58332.332
$ws.Range("J66").Value = 58332.332
$ws.Range("N66").Value = -181860.996
$ws.Range("L102").Value = 3270.7144
$ws.Range("N102").Value = -6514.7144
$ws.Range("H102").Value = 3299.524
$ws.Range("J102").Value = 3270.7144
$ws.Range("N122").Value = -20322.1432
$ws.Range("M122").Value = -7673.650000000001
$ws.Range("K122").Value = 10123.65
$ws.Range("H122").Value = 3832.4443
$ws.Range("J122").Value = 5140.7144
$ws.Range("L122").Value = 15422.1432
$ws.Range("I122").Value = 3374.55
$ws.Range("M126").Value = -14987.9228
$ws.Range("N126").Value = -26794.1
$ws.Range("I126").Value = 5819.3076
$ws.Range("J126").Value = 7284.7
$ws.Range("K126").Value = 17457.9228
$ws.Range("H126").Value = 6456.4346
$ws.Range("L126").Value = 21854.1
$ws.Range("K132").Value = 118800.552
$ws.Range("H132").Value = 25397.818
$ws.Range("I132").Value = 39600.184
$ws.Range("M132").Value = -116270.552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M16").Value = -4813.5386
$ws.Range("L16").Value = 1000
$ws.Range("H16").Value = 4699
$ws.Range("I16").Value = 4983.5386
$ws.Range("K16").Value = 4983.5386
$ws.Range("N16").Value = -1340
$ws.Range("J16").Value = 1000
$ws.Range("I22").Value = 111628.7
$ws.Range("H22").Value = 59730.473
$ws.Range("K22").Value = 111628.7
$ws.Range("M22").Value = -111333.7
$ws.Range("M27").Value = -111521.7
$ws.Range("H27").Value = 59730.473
$ws.Range("I27").Value = 111628.7
$ws.Range("K27").Value = 111628.7
$ws.Range("H68").Value = 2524.5
$ws.Range("I68").Value = 0
$ws.Range("N68").Value = -4022.5
$ws.Range("M68").ClearContents()
$ws.Range("J68").Value = 2524.5
$ws.Range("L68").Value = 2524.5
$ws.Range("K68").Value = 0
$ws.Range("H71").Value = 2524.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 12622.5
$ws.Range("I71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("J71").Value = 2524.5
$ws.Range("N71").Value = -20110.5
$ws.Range("K132").Value = 185644.569
$ws.Range("H132").Value = 59432.363
$ws.Range("I132").Value = 61881.523
$ws.Range("M132").Value = -183114.569
$ws.Range("I136").Value = 1782.6842
$ws.Range("M136").Value = -2798.0526
$ws.Range("H136").Value = 1883.55
$ws.Range("K136").Value = 5348.0526

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 8977
$ws.Range("L5").Value = 9721.25
$ws.Range("N5").Value = -9945.25
$ws.Range("J5").Value = 9721.25
$ws.Range("N7").Value = -5725.5
$ws.Range("K7").Value = 8000
$ws.Range("J7").Value = 5499.5
$ws.Range("I7").Value = 8000
$ws.Range("M7").Value = -7887
$ws.Range("L7").Value = 5499.5
$ws.Range("H7").Value = 6333
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("L12").Value = 0
$ws.Range("I51").Value = 5499
$ws.Range("J51").Value = 29900
$ws.Range("M51").Value = -4989
$ws.Range("L51").Value = 29900
$ws.Range("N51").Value = -30920
$ws.Range("K51").Value = 5499
$ws.Range("H51").Value = 17699.5
$ws.Range("N62").Value = -136532.62
$ws.Range("L62").Value = 135284.62
$ws.Range("J62").Value = 135284.62
$ws.Range("H62").Value = 85497.16
$ws.Range("M62").Value = -5213.2
$ws.Range("I62").Value = 5837.2
$ws.Range("K62").Value = 5837.2
$ws.Range("L63").Value = 28332.166
$ws.Range("H63").Value = 28332.166
$ws.Range("J63").Value = 28332.166
$ws.Range("N63").Value = -29580.166
$ws.Range("K65").Value = 29186
$ws.Range("I65").Value = 5837.2
$ws.Range("H65").Value = 85497.16
$ws.Range("L65").Value = 676423.1
$ws.Range("M65").Value = -26066
$ws.Range("J65").Value = 135284.62
$ws.Range("N65").Value = -682663.1
$ws.Range("L66").Value = 84996.49800000001
$ws.Range("H66").Value = 28332.166
$ws.Range("J66").Value = 28332.166
$ws.Range("N66").Value = -91236.49800000001
$ws.Range("J107").Value = 2112.5
$ws.Range("I107").Value = 1400
$ws.Range("L107").Value = 6337.5
$ws.Range("M107").Value = -2280
$ws.Range("H107").Value = 1756.25
$ws.Range("N107").Value = -10177.5
$ws.Range("K107").Value = 4200
$ws.Range("N122").Value = -22727.5
$ws.Range("M122").Value = -4114.875100000001
$ws.Range("K122").Value = 6564.875100000001
$ws.Range("H122").Value = 2477.077
$ws.Range("J122").Value = 5942.5
$ws.Range("L122").Value = 17827.5
$ws.Range("I122").Value = 2188.2917
$ws.Range("M126").Value = -83117.39
$ws.Range("I126").Value = 28529.13
$ws.Range("K126").Value = 85587.39
$ws.Range("H126").Value = 23638.17
$ws.Range("K132").Value = 50644.071
$ws.Range("H132").Value = 16467.729
$ws.Range("I132").Value = 16881.357
$ws.Range("M132").Value = -48114.071
$ws.Range("I136").Value = 3778.7827
$ws.Range("M136").Value = -8786.348100000001
$ws.Range("H136").Value = 3778.7827
$ws.Range("K136").Value = 11336.3481
